# Add "2022-Q3" sheet (fund-holdings detail) right after "总计", and
# update the "总计" summary sheet with the new quarter's row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert 2022-Q3 as the new first
#    data row and push every following quarter down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Grow the used range down to row 9 while inheriting column-A's style,
# then fill in every row (index, quarter label, count, value) bottom-up
# is unnecessary since we just overwrite all values outright.
$summary.Range("A8").Copy($summary.Range("A9"))

$summaryRows = @(
    @(0, "2022-Q3", 12, 1.09),
    @(1, "2022-Q2", 19, 1.79),
    @(2, "2022-Q1", 21, 2.65),
    @(3, "2021-Q4", 20, 6.13),
    @(4, "2021-Q3", 21, 9.83),
    @(5, "2021-Q2", 26, 14.91),
    @(6, "2021-Q1", 43, 23.18),
    @(7, "2020-Q4", 27, 16.29)
)

$r = 2
foreach ($row in $summaryRows) {
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" detail sheet. Clone the "2022-Q2" sheet so
#    the header row / column formatting / page setup match the other
#    quarterly sheets exactly, then overwrite its data with the new
#    quarter's fund-holding rows.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item(2)
$template.Copy($null, $summary)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Drop the old sheet's extra rows (it had 19 data rows, the new one only
# has 12) so the dimension / row count match the new data.
$q3.Range("A14:H20").Delete()

# Fund code / size / position columns are text in the source data (they
# keep trailing zeros & leading zeros), so mark them as Text before
# writing. The rank columns (A, H) stay numeric.
$q3.Range("B2:G13").NumberFormat = "@"

$q3Rows = @(
    @(0, "160527", "博时研究优选3年封闭运作灵活配置混合A", "16.93", "97.11", "5.83", "0.9870", 6),
    @(1, "160528", "博时研究优选3年封闭运作灵活配置混合C", "0.59", "97.11", "5.83", "0.0344", 6),
    @(2, "004223", "金信多策略精选灵活配置混合", "0.32", "92.79", "4.76", "0.0152", 10),
    @(3, "013733", "红塔红土盛丰混合A", "0.37", "68.35", "3.66", "0.0135", 4),
    @(4, "000743", "红塔红土盛世普益灵活配置混合", "1.12", "20.22", "1.07", "0.0120", 8),
    @(5, "010663", "长江均衡成长混合A", "0.20", "81.80", "4.74", "0.0095", 1),
    @(6, "003659", "山西证券策略精选灵活配置混合", "0.27", "78.35", "2.69", "0.0073", 7),
    @(7, "013734", "红塔红土盛丰混合C", "0.10", "68.35", "3.66", "0.0037", 4),
    @(8, "002023", "红塔红土稳健回报灵活配置混合A", "0.09", "62.51", "3.40", "0.0031", 7),
    @(9, "010664", "长江均衡成长混合C", "0.05", "81.80", "4.74", "0.0024", 1),
    @(10, "004696", "东兴量化优享灵活配置混合", "0.02", "93.60", "3.93", "0.0008", 6),
    @(11, "002024", "红塔红土稳健回报灵活配置混合C", "0.00", "62.51", "3.40", 0, 7)
)

$r = 2
foreach ($row in $q3Rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Row 13's "持有市值" happens to be exactly 0 in the source data, stored
# as a plain number rather than the text "0.0000" used elsewhere in that
# column - match that quirk exactly.
$q3.Cells.Item(13, 7).Value = 0

# Restore the original active sheet/selection.
$summary.Activate()
$summary.Range("A1").Select()
